# Stock update by raj, time 12:13
# Updates Quantity (B), Rate (C), % (D) and Value (E) columns for a set of
# stock-item rows on the "PATRIKA 24-25" sheet, and refreshes the grand
# totals row (638) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PATRIKA 24-25")

# Row 69
$ws.Cells.Item(69, 2).Value = 124
$ws.Cells.Item(69, 3).Value = 197.5
$ws.Cells.Item(69, 5).Value = 138.25

# Row 78
$ws.Cells.Item(78, 2).Value = 148
$ws.Cells.Item(78, 3).Value = 267.5

# Row 116
$ws.Cells.Item(116, 2).Value = 31
$ws.Cells.Item(116, 3).Value = -1.5
$ws.Cells.Item(116, 5).Value = -18.75

# Row 162
$ws.Cells.Item(162, 2).Value = 24
$ws.Cells.Item(162, 3).Value = 54.5
$ws.Cells.Item(162, 5).Value = 152.6

# Row 167
$ws.Cells.Item(167, 2).Value = 56
$ws.Cells.Item(167, 3).Value = 45
$ws.Cells.Item(167, 5).Value = 135

# Row 189
$ws.Cells.Item(189, 2).Value = 100
$ws.Cells.Item(189, 3).Value = 18.5
$ws.Cells.Item(189, 5).Value = 83.25

# Row 190
$ws.Cells.Item(190, 2).Value = 86
$ws.Cells.Item(190, 3).Value = 27
$ws.Cells.Item(190, 5).Value = 101.86

# Row 213
$ws.Cells.Item(213, 2).Value = 120
$ws.Cells.Item(213, 3).Value = 1.5
$ws.Cells.Item(213, 5).Value = 6.17

# Row 227
$ws.Cells.Item(227, 2).Value = 112
$ws.Cells.Item(227, 3).Value = -6
$ws.Cells.Item(227, 5).Value = -27

# Row 237
$ws.Cells.Item(237, 2).Value = 106
$ws.Cells.Item(237, 3).Value = 37
$ws.Cells.Item(237, 5).Value = 222

# Row 293
$ws.Cells.Item(293, 2).Value = 125
$ws.Cells.Item(293, 3).Value = 6.5
$ws.Cells.Item(293, 5).Value = 40.17

# Row 334
$ws.Cells.Item(334, 2).Value = 56
$ws.Cells.Item(334, 3).Value = 17
$ws.Cells.Item(334, 5).Value = 117.48

# Row 377
$ws.Cells.Item(377, 2).Value = 14
$ws.Cells.Item(377, 3).Value = 1
$ws.Cells.Item(377, 5).Value = 10.5

# Row 424
$ws.Cells.Item(424, 2).Value = 6
$ws.Cells.Item(424, 3).Value = 25
$ws.Cells.Item(424, 5).Value = 95

# Row 487
$ws.Cells.Item(487, 2).Value = 66
$ws.Cells.Item(487, 3).Value = -9.5
$ws.Cells.Item(487, 5).Value = -90.25

# Row 513
$ws.Cells.Item(513, 2).Value = 118
$ws.Cells.Item(513, 3).Value = 47
$ws.Cells.Item(513, 5).Value = 107.16

# Row 537
$ws.Cells.Item(537, 2).Value = 90
$ws.Cells.Item(537, 3).Value = 17.5
$ws.Cells.Item(537, 5).Value = 63

# Row 561 - previously had only Quantity filled in; now Rate/%/Value populate too.
# Copy the Rate/%/Value number-formats from the row above (same pattern as
# every other populated row) before writing the new figures.
$ws.Range("C560:E560").Copy()
$ws.Range("C561:E561").PasteSpecial(-4122)
$ws.Cells.Item(561, 2).Value = 76
$ws.Cells.Item(561, 3).Value = -7
$ws.Cells.Item(561, 4).Value = 5.5
$ws.Cells.Item(561, 5).Value = -38.5

# Row 589
$ws.Cells.Item(589, 2).Value = 7
$ws.Cells.Item(589, 3).Value = 13.5
$ws.Cells.Item(589, 5).Value = 40.5

# Row 613
$ws.Cells.Item(613, 2).Value = 71
$ws.Cells.Item(613, 3).Value = 12.3

# Grand total row
$ws.Cells.Item(638, 3).Value = 37557.51
$ws.Cells.Item(638, 5).Value = 67564.7
